$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030726297540218
$ws.Range("D2").Value = 1.033966068688836
$ws.Range("E2").Value = 1.030315697076221
$ws.Range("F2").Value = 1.036135004436687
$ws.Range("I2").Value = 1.032339117213274
$ws.Range("J2").Value = 1.035865950840768
$ws.Range("K2").Value = 1.036766899286371
$ws.Range("L2").Value = 1.033127063265789
$ws.Range("M2").Value = 1.038929613308421
$ws.Range("N2").Value = 1.037336999306461

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032294342180037
$ws.Range("D3").Value = 1.035163171379579
$ws.Range("E3").Value = 1.031669516719321
$ws.Range("F3").Value = 1.038474852550451
$ws.Range("I3").Value = 1.03274007011188
$ws.Range("J3").Value = 1.037072418446276
$ws.Range("K3").Value = 1.03777205743496
$ws.Range("L3").Value = 1.034287746945129
$ws.Range("M3").Value = 1.041074943880548
$ws.Range("N3").Value = 1.038545180234339

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033306671293695
$ws.Range("D4").Value = 1.035935506472467
$ws.Range("E4").Value = 1.032543683083307
$ws.Range("F4").Value = 1.039981842094572
$ws.Range("I4").Value = 1.032996811771873
$ws.Range("J4").Value = 1.037850468466743
$ws.Range("K4").Value = 1.038419634278845
$ws.Range("L4").Value = 1.035036419991269
$ws.Range("M4").Value = 1.042455778652844
$ws.Range("N4").Value = 1.039324335175078

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033731715264339
$ws.Range("D5").Value = 1.03625966087158
$ws.Range("E5").Value = 1.032910749803587
$ws.Range("F5").Value = 1.040613733894353
$ws.Range("I5").Value = 1.033104103803129
$ws.Range("J5").Value = 1.038176943605635
$ws.Range("K5").Value = 1.03869120642058
$ws.Range("L5").Value = 1.035350603840775
$ws.Range("M5").Value = 1.043034563216151
$ws.Range("N5").Value = 1.039651273946114

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033803050733513
$ws.Range("D6").Value = 1.036314056671154
$ws.Range("E6").Value = 1.032972356799872
$ws.Range("F6").Value = 1.040719735925752
$ws.Range("I6").Value = 1.033122081083718
$ws.Range("J6").Value = 1.038231724291079
$ws.Range("K6").Value = 1.038736765562787
$ws.Range("L6").Value = 1.03540332419342
$ws.Range("M6").Value = 1.043131644025221
$ws.Range("N6").Value = 1.039706132426416

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033312352854729
$ws.Range("D7").Value = 1.035939839931243
$ws.Range("E7").Value = 1.032548589535038
$ws.Range("F7").Value = 1.039990291879062
$ws.Range("I7").Value = 1.032998247930161
$ws.Range("J7").Value = 1.037854833255858
$ws.Range("K7").Value = 1.038423265657691
$ws.Range("L7").Value = 1.035040620313448
$ws.Range("M7").Value = 1.04246351910132
$ws.Range("N7").Value = 1.039328706162694

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031256708931827
$ws.Range("D8").Value = 1.034371110249953
$ws.Range("E8").Value = 1.030773614319127
$ws.Range("F8").Value = 1.036927251300846
$ws.Range("I8").Value = 1.032475183324518
$ws.Range("J8").Value = 1.036274229151046
$ws.Range("K8").Value = 1.037107187397051
$ws.Range("L8").Value = 1.03351981688095
$ws.Range("M8").Value = 1.039656180303072
$ws.Range("N8").Value = 1.037745857418757

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027616258844152
$ws.Range("D9").Value = 1.031589045762062
$ws.Range("E9").Value = 1.027631340902311
$ws.Range("F9").Value = 1.031473944382903
$ws.Range("I9").Value = 1.031532575307965
$ws.Range("J9").Value = 1.033468567068438
$ws.Range("K9").Value = 1.034766080100155
$ws.Range("L9").Value = 1.030821462193771
$ws.Range("M9").Value = 1.034651357771974
$ws.Range("N9").Value = 1.034936210974177

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025176357823643
$ws.Range("D10").Value = 1.029721865727992
$ws.Range("E10").Value = 1.025526143103326
$ws.Range("F10").Value = 1.027798245861186
$ws.Range("I10").Value = 1.030889820239087
$ws.Range("J10").Value = 1.031583812668208
$ws.Range("K10").Value = 1.033190046329325
$ws.Range("L10").Value = 1.029009585721149
$ws.Range("M10").Value = 1.031273406576361
$ws.Range("N10").Value = 1.033048780006514

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024116620663884
$ws.Range("D11").Value = 1.028910279322375
$ws.Range("E11").Value = 1.024611985738713
$ws.Range("F11").Value = 1.026196512039064
$ws.Range("I11").Value = 1.030608029516095
$ws.Range("J11").Value = 1.030764169361866
$ws.Range("K11").Value = 1.032503866186589
$ws.Range("L11").Value = 1.028221826409371
$ws.Range("M11").Value = 1.029800343335674
$ws.Range("N11").Value = 1.032227972712662

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023722485374168
$ws.Range("D12").Value = 1.028608346096102
$ws.Range("E12").Value = 1.024272026221178
$ws.Range("F12").Value = 1.02559998323757
$ws.Range("I12").Value = 1.030502832149464
$ws.Range("J12").Value = 1.030459175422984
$ws.Range("K12").Value = 1.032248415563668
$ws.Range("L12").Value = 1.027928725597951
$ws.Range("M12").Value = 1.029251573289651
$ws.Range("N12").Value = 1.031922545647404

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023807051668285
$ws.Range("D13").Value = 1.028673133411444
$ws.Range("E13").Value = 1.024344967009181
$ws.Range("F13").Value = 1.02572801273064
$ws.Range("I13").Value = 1.030525421318817
$ws.Range("J13").Value = 1.03052462236337
$ws.Range("K13").Value = 1.032303236667151
$ws.Range("L13").Value = 1.027991619146767
$ws.Range("M13").Value = 1.029369359909168
$ws.Range("N13").Value = 1.031988085529953

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024084051636627
$ws.Range("D14").Value = 1.028885331155761
$ws.Range("E14").Value = 1.024583892819993
$ws.Range("F14").Value = 1.026147235182429
$ws.Range("I14").Value = 1.030599344669198
$ws.Range("J14").Value = 1.030738969588338
$ws.Range("K14").Value = 1.03248276229987
$ws.Range("L14").Value = 1.028197608700923
$ws.Range("M14").Value = 1.02975501493239
$ws.Range("N14").Value = 1.032202737152566

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024254653543255
$ws.Range("D15").Value = 1.029016010079717
$ws.Range("E15").Value = 1.024731049453344
$ws.Range("F15").Value = 1.026405321909421
$ws.Range("I15").Value = 1.030644821176007
$ws.Range("J15").Value = 1.030870963839291
$ws.Range("K15").Value = 1.03259329777459
$ws.Range("L15").Value = 1.028324460163899
$ws.Range("M15").Value = 1.029992415088466
$ws.Range("N15").Value = 1.032334918850494

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025246619521401
$ws.Range("D16").Value = 1.029775662160909
$ws.Range("E16").Value = 1.025586757134148
$ws.Range("F16").Value = 1.027904329830811
$ws.Range("I16").Value = 1.030908448077719
$ws.Range("J16").Value = 1.031638134318934
$ws.Range("K16").Value = 1.033235506064406
$ws.Range("L16").Value = 1.029061798302232
$ws.Range("M16").Value = 1.031370945984456
$ws.Range("N16").Value = 1.033103178800215

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025867975000432
$ws.Range("D17").Value = 1.030251338573239
$ws.Range("E17").Value = 1.026122817901934
$ws.Range("F17").Value = 1.02884186933016
$ws.Range("I17").Value = 1.031072880179032
$ws.Range("J17").Value = 1.032118406526788
$ws.Range("K17").Value = 1.033637336269803
$ws.Range("L17").Value = 1.029523445684097
$ws.Range("M17").Value = 1.032232847157887
$ws.Range("N17").Value = 1.033584133049692

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026230089127442
$ws.Range("D18").Value = 1.030528495616554
$ws.Range("E18").Value = 1.026435243674589
$ws.Range("F18").Value = 1.029387746170924
$ws.Range("I18").Value = 1.03116845589857
$ws.Range("J18").Value = 1.032398201140317
$ws.Range("K18").Value = 1.033871356157157
$ws.Range("L18").Value = 1.029792408283091
$ws.Range("M18").Value = 1.03273458052382
$ws.Range("N18").Value = 1.033864325003665

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026353508196021
$ws.Range("D19").Value = 1.030622948948397
$ws.Range("E19").Value = 1.026541730823319
$ws.Range("F19").Value = 1.029573712383158
$ws.Range("I19").Value = 1.031200988147327
$ws.Range("J19").Value = 1.032493546534147
$ws.Range("K19").Value = 1.033951089975153
$ws.Range("L19").Value = 1.029884065563867
$ws.Range("M19").Value = 1.032905490628323
$ws.Range("N19").Value = 1.033959805798891

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02580134181444
$ws.Range("D20").Value = 1.03020033378693
$ws.Range("E20").Value = 1.026065329532821
$ws.Range("F20").Value = 1.028741381337454
$ws.Range("I20").Value = 1.031055272837177
$ws.Range("J20").Value = 1.032066913117187
$ws.Range("K20").Value = 1.033594261089156
$ws.Range("L20").Value = 1.029473947274803
$ws.Range("M20").Value = 1.032140476996024
$ws.Range("N20").Value = 1.033532566513543

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024002496016257
$ws.Range("D21").Value = 1.028822857375088
$ws.Range("E21").Value = 1.024513546277868
$ws.Range("F21").Value = 1.02602382847361
$ws.Range("I21").Value = 1.030577590701106
$ws.Range("J21").Value = 1.030675864723755
$ws.Range("K21").Value = 1.032429912366148
$ws.Range("L21").Value = 1.028136963575504
$ws.Range("M21").Value = 1.029641493940848
$ws.Range("N21").Value = 1.032139542671837

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022868580202765
$ws.Range("D22").Value = 1.027954036192125
$ws.Range("E22").Value = 1.02353555415849
$ws.Range("F22").Value = 1.024306064960743
$ws.Range("I22").Value = 1.030274196624589
$ws.Range("J22").Value = 1.029798116324696
$ws.Range("K22").Value = 1.031694521307858
$ws.Range("L22").Value = 1.027293497444247
$ws.Range("M22").Value = 1.028060951770725
$ws.Range("N22").Value = 1.031260547769372

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02346997111175
$ws.Range("D23").Value = 1.028414878640311
$ws.Range("E23").Value = 1.024054230397335
$ws.Range("F23").Value = 1.025217566723101
$ws.Range("I23").Value = 1.030435323262729
$ws.Range("J23").Value = 1.030263728888343
$ws.Range("K23").Value = 1.032084683915032
$ws.Range("L23").Value = 1.02774090860074
$ws.Range("M23").Value = 1.028899727971999
$ws.Range("N23").Value = 1.031726821556266

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025831451447001
$ws.Range("D24").Value = 1.030223381571971
$ws.Range("E24").Value = 1.02609130682291
$ws.Range("F24").Value = 1.028786790543523
$ws.Range("I24").Value = 1.031063229871298
$ws.Range("J24").Value = 1.03209018182313
$ws.Range("K24").Value = 1.03361372602417
$ws.Range("L24").Value = 1.029496314428576
$ws.Range("M24").Value = 1.032182218180908
$ws.Range("N24").Value = 1.033555868263718

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028559627718168
$ws.Range("D25").Value = 1.032310434959841
$ws.Range("E25").Value = 1.028445475527907
$ws.Range("F25").Value = 1.032890638687031
$ws.Range("I25").Value = 1.031778769256901
$ws.Range("J25").Value = 1.034196378262764
$ws.Range("K25").Value = 1.035373971489905
$ws.Range("L25").Value = 1.031521298798713
$ws.Range("M25").Value = 1.035952339425175
$ws.Range("N25").Value = 1.03566505574388
